$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E14").Value = 15944
